# Kế hoạch làm việc
#
# The "Launch" task (row 13) status cell F13 still said "Not Started" (red
# text on a colored fill) while every other task in the STATUS column
# (F3:F12) already says "Completed" (accent-colored text, no fill).
# Bring F13 in line with the rest of the column: copy the existing
# "Completed" formatting from F3 onto F13 and set its text to "Completed".
# Finally, extend the sheet's remembered selection down to F13 so it covers
# the whole STATUS column again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("F3")
$target = $ws.Range("F13")

# Reuse the same cell formatting ("Completed" style) already used by F3:F12.
$source.Copy() | Out-Null
$target.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$target.Value = "Completed"

# Restore the STATUS column selection to include the now-updated F13.
$ws.Range("F9:F13").Select() | Out-Null
